$p = $ppt.ActivePresentation

# --- Add a new slide at the end (slide 9), same "Blank" layout used by the
# rest of this deck's slides -------------------------------------------------
$layout = $p.Slides.Item(1).CustomLayout
$s = $p.Slides.Add($p.Slides.Count + 1, 12)
$s.CustomLayout = $layout

# --- Add the "CaixaDeTexto 1" textbox with the DATEDIFF/DATE_ADD/DATE_SUB
# explanations, matching the style used by the other "calculos e funcoes"
# slides in this deck ---------------------------------------------------------
$tb = $s.Shapes.AddTextbox(1, 224.85354330708662, 201.96094488188976, 694.5653543307087, 123.59527559055118)
$tb.Name = "CaixaDeTexto 1"
$tb.Fill.Visible = $false

$tf = $tb.TextFrame
$tf.WordWrap = $true
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.LanguageID = "pt-BR"
$tr.Text = "DATEDIFF() - RETORNA A QUANTIDADE DE DIAS ENTRE DUAS DATAS`rDATE_ADD() - ADICIONA`rDATE_SUB() -"
$tr.Font.Size = 24

$para1 = $tr.Paragraphs(1, 1)
$para1.Characters(1, 13).Font.Bold = $true

$para2 = $tr.Paragraphs(2, 1)
$para2.Characters(1, 13).Font.Bold = $true

$para3 = $tr.Paragraphs(3, 1)
$para3.Characters(1, 12).Font.Bold = $true

# --- Give the new slide an (empty) notes page, mirroring the other slides
# in the deck, which all carry an (empty) notes page -------------------------
$np = $s.NotesPage
$body = $np.Shapes.Placeholders.Item(2)
$body.TextFrame.TextRange.Text = ""
